# Re-generated quadratic/linear problem data (new random draw).
# Mirrors the commit "volver a generar problemas cuadraticos y lineales":
# the numeric parameters on several sheets are replaced with freshly
# generated values. The expression cells in column A already contain
# non-numeric text (e.g. "-12.28... - 2x_1 + ..."), so Excel stores them
# as text automatically; the purely numeric-looking cells need to be
# forced to Text format first so they keep being stored as shared
# strings (matching the source workbook, where every such cell is
# t="s") instead of being auto-coerced into numbers.

$wb = $excel.ActiveWorkbook

# ---- Sheet: Restricciones_del_follower ----
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

# Force the purely-numeric-looking cells to stay text, like the source file.
$ws.Range("B2:B4").NumberFormat = "@"
$ws.Range("D2:F4").NumberFormat = "@"

$ws.Range("A2").Value = "-12.283442485306468 - 2x_1 + 3.1187237615449206y_1 + 1.9204030226700257y_2"
$ws.Range("B2").Value = "14.783442485306468"
$ws.Range("D2").Value = "0.96"
$ws.Range("F2").Value = "6.4"

$ws.Range("A3").Value = "5.531040736639048 + x_1 - 3x_2 - 0.8984582512388339y_1 - 0.23841542738325772y_2"
$ws.Range("B3").Value = "-7.531040736639048"
$ws.Range("D3").Value = "0.9"
$ws.Range("E3").Value = "4.2"

$ws.Range("A4").Value = "-5.876876574307305 + x_1 + x_2 + 0.2795969773299748y_1 + 0.3853904282115869y_2"
$ws.Range("B4").Value = "3.8168765743073045"
$ws.Range("D4").Value = "0.28"
$ws.Range("E4").Value = "8.2"
$ws.Range("F4").Value = "0.8999999999999999"

# ---- Sheet: Punto_modificado ----
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2:D2").NumberFormat = "@"
$ws.Range("A2").Value = "2.55"
$ws.Range("B2").Value = "1.25"
$ws.Range("C2").Value = "4.25"
$ws.Range("D2").Value = "2.15"

# ---- Sheet: Vector_bf (index 5). ----
# NOTE: "Vector_bf" and "Vector_BF" differ only by case, and
# Worksheets.Item(name) resolves case-insensitively to the first match
# ("Vector_bf"), so the two sheets must be addressed by 1-based index.
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("A2").Value = "1.736350461379434"
$ws.Range("A3").Value = "-2.7369223370175373"

# ---- Sheet: Vector_BF (index 6) ----
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2:A5").NumberFormat = "@"
$ws.Range("A2").Value = "-10.399999999999999"
$ws.Range("A3").Value = "3.400000000000002"
$ws.Range("A4").Value = "0.9808294410973093"
$ws.Range("A5").Value = "-2.15885671632533"

# ---- Sheet: Vector_Alpha ----
# A2/A3 on this sheet are genuine numbers in the source (no t="s"), so
# they are written as real numeric values.
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 1.1099999999999999
$ws.Range("A3").Value = 1.53
